$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

# Change TradeLink column (H) values from "U" to "B" for rows 4-11
for ($r = 4; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Value() -eq "U") {
        $cell.Value = "B"
    }
}

# Update the active selection on the sheet to match the recorded view state
$ws.Activate()
$ws.Range("M9").Select()
